$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-04-17"

# Update header label for column B (April 2022 through date)
$ws.Range("B1").Value = "April 2022 (through April 17)"

# Apply per-cell updates (row, cell, new value)
$updates = @{
    "F2"  = 6
    "J2"  = 6
    "V2"  = 3
    "AD2" = 3
    "B3"  = 5
    "J3"  = 1
    "Z3"  = 1
    "B4"  = 3
    "F4"  = 5
    "F5"  = 3
    "V5"  = 4
    "N8"  = 2
    "F11" = 3
    "F15" = 2
    "Z17" = 1
    "J23" = 1
    "B44" = 2
    "V50" = 1
    "B53" = 2
    "R54" = 1
    "F60" = 2
    "D61" = 1
    "B73" = 1
    "B91" = 1
    "J93" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
